# Apply the "break out stock.yaml completed" edit:
#  1. Q59 / Q63 / Q67 "detect_structure" values reset to 0.
#  2. R1485 / R1486 "backup" placeholders backfilled from blank to 0.
#  3. 16 new weekly rows (1487-1502) appended, continuing the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Reset stray detect_structure (Q) values back to 0 ---
$ws.Cells.Item(59, 17).Value = 0
$ws.Cells.Item(63, 17).Value = 0
$ws.Cells.Item(67, 17).Value = 0

# --- 2. Backfill the last two existing rows' backup (R) column with 0 ---
$ws.Cells.Item(1485, 18).Value = 0
$ws.Cells.Item(1486, 18).Value = 0

# --- 3. Append the new weekly OHLCV rows ---
# Columns A..Q (R/backup is left blank, matching the freshly-appended rows'
# placeholder state):
#   A Datetime(serial) B Open C High D Low E Close F Adj Close G Volume
#   H Year I Month J Day K Hour L Minute M Second N Week O isPivot
#   P two_line_structure Q detect_structure
$newRows = @(
    @(45474, 1214.050048828125, 1311.800048828125, 1203, 1303.550048828125, 1303.550048828125, 6766511, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0),
    @(45481, 1305.099975585938, 1344.199951171875, 1285.75, 1326.849975585938, 1326.849975585938, 5459020, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0),
    @(45488, 1360, 1409.900024414062, 1320.099975585938, 1330.800048828125, 1330.800048828125, 10079759, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 0),
    @(45495, 1332, 1389, 1314.449951171875, 1386.199951171875, 1386.199951171875, 3950950, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(45502, 1390, 1450.800048828125, 1387.599975585938, 1443.300048828125, 1443.300048828125, 5871152, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
    @(45509, 1428.349975585938, 1488.300048828125, 1395, 1449.699951171875, 1449.699951171875, 9482969, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0),
    @(45516, 1460, 1533.199951171875, 1421.449951171875, 1502.75, 1502.75, 12746459, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0),
    @(45523, 1502.099975585938, 1547, 1490.099975585938, 1538.199951171875, 1538.199951171875, 6291441, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0),
    @(45530, 1538, 1580.599975585938, 1531.099975585938, 1569.400024414062, 1569.400024414062, 4502809, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0),
    @(45537, 1580.800048828125, 1584.449951171875, 1516.25, 1537.400024414062, 1537.400024414062, 4899562, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0),
    @(45544, 1537.349975585938, 1592, 1502.400024414062, 1567.25, 1567.25, 7951867, 2024, 9, 9, 0, 0, 0, 37, 1, 0, 0),
    @(45551, 1574.849975585938, 1583.949951171875, 1482.599975585938, 1496.349975585938, 1496.349975585938, 5144655, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0),
    @(45558, 1498, 1531.449951171875, 1470, 1511.349975585938, 1511.349975585938, 6096222, 2024, 9, 23, 0, 0, 0, 39, 0, 0, 0),
    @(45565, 1512.050048828125, 1518.900024414062, 1407, 1466.050048828125, 1466.050048828125, 6910700, 2024, 9, 30, 0, 0, 0, 40, 0, 0, 0),
    @(45572, 1467, 1520.050048828125, 1439.5, 1495.849975585938, 1495.849975585938, 4437619, 2024, 10, 7, 0, 0, 0, 41, 0, 0, 0),
    @(45579, 1496.099975585938, 1502.699951171875, 1450, 1472.699951171875, 1472.699951171875, 3096009, 2024, 10, 14, 0, 0, 0, 42, 0, 0, 0)
)

$r = 1487
foreach ($row in $newRows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    # Column A carries the same date/time display format as the rest of the sheet.
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $r = $r + 1
}
